# Adds a "standard error of the mean" (FeltLocation_sem) column next to the
# existing FeltLocation column on each of the per-duration sheets, renaming
# the existing FeltLocation column to FeltLocation_mean.
#
#  - Duration_0.1 (sheet index 2): D=FeltLocation -> FeltLocation_mean, new E=FeltLocation_sem
#  - Duration_1.0 (sheet index 3): D=FeltLocation -> FeltLocation_mean, new E=FeltLocation_sem
#  - Duration_2.0 (sheet index 4): D=FeltLocation -> FeltLocation_mean, new E=FeltLocation_sem
#  - AllDurations (sheet index 5): C=FeltLocation -> FeltLocation_mean, new D=FeltLocation_sem
#    (AllDurations has no Duration column, so the new sem column lands one
#    letter earlier than on the per-duration sheets.)

$wb = $excel.ActiveWorkbook

$semBySheet = @{
  "Duration_0.1" = @(0.110171315, 0.03908652154687299, 0.08103833338678147, 0.2157364296148972, 0.06591565264255492, 0.05294134990086616, 0.1951076997441413, 0.1358648912345851, 0.103304458405211, 0.08770418712734175, 0.03491347755397142, 0.04029213953425036, 0.004543166666666654, 0.125, 0.004038924999999998)
  "Duration_1.0" = @(0.07350219628285615, 0.0253547637660217, 0.04399245149670595, 0.08140434709843383, 0.0938993718398614, 0.0613428999862325, 0.1206212465522831, 0.1615808256107296, 0.1305259416407927, 0.02919907643785944, 0.02204273727897695, 0.02526908665852078, 0.02376171338888728, 0.03774461303435674, 0.02728231192520089)
  "Duration_2.0" = @(0.01725933582406557, 0.02915320620593256, 0.04017178979663463, 0.06011087313284629, 0.02451828750385665, 0.0966035480163029, 0.1057900180329608, 0.07145830211582019, 0.1432080587623498, 0.06511646377055884, 0.05550644451455427, 0.1736949542517715, 0.046062875, 0.03463314494889013, 0.01823615000000001)
  "AllDurations" = @(0.03105509752217602, 0.01791172485067583, 0.0325585621551735, 0.07841287676813033, 0.04022148535914662, 0.05164727619696929, 0.07860977714961188, 0.06836843191938237, 0.06671312960031939, 0.03338423279515879, 0.02805226666324806, 0.06318681993151082, 0.01768305799321687, 0.04122894521269786, 0.01253992107185321)
}

# column that currently holds "FeltLocation" (mean) on each sheet -- 4 (D) for
# the per-duration sheets (which also have a Duration column in C), 3 (C) for
# the combined AllDurations sheet.
$meanColBySheet = @{
  "Duration_0.1" = 4
  "Duration_1.0" = 4
  "Duration_2.0" = 4
  "AllDurations" = 3
}

foreach ($sheetName in @("Duration_0.1", "Duration_1.0", "Duration_2.0", "AllDurations")) {
  $ws = $wb.Worksheets.Item($sheetName)

  $meanCol = $meanColBySheet[$sheetName]
  $semCol = $meanCol + 1

  # Rename the existing header and add the new "_sem" header next to it,
  # copying the bold/centered/bordered header formatting across.
  $ws.Cells.Item(1, $meanCol).Value = "FeltLocation_mean"

  $ws.Cells.Item(1, $meanCol).Copy() | Out-Null
  $ws.Cells.Item(1, $semCol).PasteSpecial(-4122) | Out-Null
  $ws.Cells.Item(1, $semCol).Value = "FeltLocation_sem"
  $excel.CutCopyMode = 0

  $sem = $semBySheet[$sheetName]
  for ($i = 0; $i -lt $sem.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, $semCol).Value = $sem[$i]
  }
}
